# edit.ps1
# Applies the "Added BLE Library" edit to instructions.docx:
#   1. Splits the final "npm install ... --save" run into three runs
#      (space-separated "...picker ", an en dash "–", and "save") to
#      mirror an autocorrect-style split of "--save" into an en-dash + text.
#   2. Inserts a new "Install react-native-ble-plx ..." step (with its
#      yarn/npm command paragraphs) right after that line.
#   3. Moves the existing "Link all native dependencies..." /
#      "react-native link react-native-gesture-handler" paragraphs to
#      after the new BLE install step.
#   4. Appends a new final paragraph "react-native link react-native-ble-plx"
#      that now carries the trailing _GoBack bookmark.

$d = $word.ActiveDocument

# Locate the two anchor paragraphs by their (unique) text rather than a
# hard-coded index, so the script is resilient to minor paragraph shifts.
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "npm install react-native-modal-datetime-picker --save*") {
        $startPara = $p
    }
    if ($startPara -ne $null -and $t -like "react-native link react-native-gesture-handler*") {
        $endPara = $p
        break
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate the anchor paragraphs for the BLE-library edit."
}

$target = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t xml:space="preserve">npm install react-native-modal-datetime-picker </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:t>save</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Install react-native-ble-plx to the project in the “Pillbox” folder with the following command:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>yarn add react-naïve-ble-plx</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>or</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>npm install --save react-native-ble-plx</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Link all native dependencies by running the following command in the “Pillbox” folder:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>react-native link react-native-gesture-handler</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>react-native link react-native-ble-plx</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$null = $target.InsertXML($xml)
Write-Host "Applied BLE library edit."
